# Update "distance from Dma50" sheet data (2025-10-30 11:48 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("distance from Dma50")

# Rows 5/6 swap stock names (NIFTYCOMMODITIES <-> CNXREALTY)
$ws.Range("B5").Value = "NIFTYCOMMODITIES"
$ws.Range("B6").Value = "CNXREALTY"

# Rows 22/23 swap stock names (NIFTYHEALTHCARE <-> CNXIT)
$ws.Range("B22").Value = "CNXIT"
$ws.Range("B23").Value = "NIFTYHEALTHCARE"

# Column C (Distance From Sma50) value updates for rows 2-30
$ws.Range("C2").Value = 9.9245
$ws.Range("C3").Value = 7.5614
$ws.Range("C4").Value = 6.6964
$ws.Range("C5").Value = 5.4651
$ws.Range("C6").Value = 5.4524
$ws.Range("C7").Value = 5.2136
$ws.Range("C8").Value = 4.4914
$ws.Range("C9").Value = 4.4253
$ws.Range("C10").Value = 3.9643
$ws.Range("C11").Value = 3.8581
$ws.Range("C12").Value = 3.4891
$ws.Range("C13").Value = 3.4834
$ws.Range("C14").Value = 3.1701
$ws.Range("C15").Value = 3.1383
$ws.Range("C16").Value = 3.0558
$ws.Range("C17").Value = 2.9132
$ws.Range("C18").Value = 2.9095
$ws.Range("C19").Value = 2.8096
$ws.Range("C20").Value = 2.4554
$ws.Range("C21").Value = 2.3947
$ws.Range("C22").Value = 1.402
$ws.Range("C23").Value = 1.3974
$ws.Range("C24").Value = 1.3361
$ws.Range("C25").Value = 1.039
$ws.Range("C26").Value = 1.0033
$ws.Range("C27").Value = 0.8643999999999999
$ws.Range("C28").Value = 0.6471
$ws.Range("C29").Value = 0.3843
$ws.Range("C30").Value = -2.1637
